$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.714.30"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "1.866.18"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "300.84"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.5314"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").Value = "0.3728"
$ws.Range("E8").Value = "  -2.17%  "
$ws.Range("D9").Value = "0.07168"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "21.50"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").Value = "0.8879"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "0.08180"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "1.872.75"
$ws.Range("E13").Value = "  +27.85%  "
$ws.Range("D14").Value = "92.60"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("D15").Value = "5.294"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "14.76"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "0.000008485"
$ws.Range("E18").Value = "  -2.13%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "26.765.72"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "4.974"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").Value = "10.61"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").Value = "6.351"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").Value = "2.298"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "145.54"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "1.731"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.01"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "113.58"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").Value = "4.682"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").Value = "4.626"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("D31").Value = "0.09120"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "0.7990"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").Value = "0.05014"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "1.172"
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("D35").Value = "2.943"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "0.6058"
$ws.Range("E36").Value = "  +4.18%  "
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "3.173"
$ws.Range("E38").Value = "  -4.81%  "
$ws.Range("D39").Value = "0.01941"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").Value = "1.064"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "6.504"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").Value = "0.5202"
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("D43").Value = "8.711"
$ws.Range("E43").Value = "  -5.14%  "
$ws.Range("D44").Value = "114.45"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("D45").Value = "0.1490"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.946"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.635"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "37.43"
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("D50").Value = "0.06056"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "61.95"
$ws.Range("E51").Value = "  -3.93%  "
